# Weekly update: two new price records for Damasco (Castle Brite) are
# inserted at the top of the data table (rows 51-52), pushing the
# existing records (old rows 51-96) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 51, shifting
# every row from 51 downward by two positions (old row 51 -> row 53, ...,
# old row 96 -> row 98).
$ws.Rows("51:52").Insert()

# --- New row 51 ---
$ws.Cells.Item(51, 1).Value  = 10
$ws.Cells.Item(51, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(51, 3).Value  = "La Araucanía"
$ws.Cells.Item(51, 4).Value  = 45264
$ws.Cells.Item(51, 5).Value  = 9
$ws.Cells.Item(51, 6).Value  = "Fruta"
$ws.Cells.Item(51, 7).Value  = 100103
$ws.Cells.Item(51, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(51, 9).Value  = 100103003
$ws.Cells.Item(51, 10).Value = "Damasco"
$ws.Cells.Item(51, 11).Value = "Castle Brite"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 80
$ws.Cells.Item(51, 14).Value = 30000
$ws.Cells.Item(51, 15).Value = 30000
$ws.Cells.Item(51, 16).Value = 30000
$ws.Cells.Item(51, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(51, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(51, 19).Value = 2000
$ws.Cells.Item(51, 20).Value = 15

# --- New row 52 ---
$ws.Cells.Item(52, 1).Value  = 10
$ws.Cells.Item(52, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(52, 3).Value  = "La Araucanía"
$ws.Cells.Item(52, 4).Value  = 45264
$ws.Cells.Item(52, 5).Value  = 9
$ws.Cells.Item(52, 6).Value  = "Fruta"
$ws.Cells.Item(52, 7).Value  = 100103
$ws.Cells.Item(52, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(52, 9).Value  = 100103003
$ws.Cells.Item(52, 10).Value = "Damasco"
$ws.Cells.Item(52, 11).Value = "Castle Brite"
$ws.Cells.Item(52, 12).Value = "Primera"
$ws.Cells.Item(52, 13).Value = 120
$ws.Cells.Item(52, 14).Value = 32000
$ws.Cells.Item(52, 15).Value = 32000
$ws.Cells.Item(52, 16).Value = 32000
$ws.Cells.Item(52, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(52, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(52, 19).Value = 1778
$ws.Cells.Item(52, 20).Value = 18

# Make sure the date cells keep the same date number format used by the
# rest of column D.
$ws.Range("D51:D52").NumberFormat = $ws.Range("D53").NumberFormat
